$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken by permuting existing rows (week re-sorted).
# Columns: D (Fecha/serial), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)

$ws.Range("D2").Value = 44449
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 400

$ws.Range("D3").Value = 44418
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("P3").Value = 500

$ws.Range("D4").Value = 44446
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("P4").Value = 467

$ws.Range("D5").Value = 44474
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 333

$ws.Range("D7").Value = 44435
$ws.Range("O7").Value = "Provincia de Limarí"

$ws.Range("O8").Value = "Provincia del Elquí"

$ws.Range("D9").Value = 44376
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 18000
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 600

$ws.Range("D10").Value = 44432
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 14000
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 467

$ws.Range("D11").Value = 44460
$ws.Range("J11").Value = 45
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 13000
$ws.Range("P11").Value = 433

$ws.Range("D12").Value = 44421
$ws.Range("J12").Value = 25
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15400
$ws.Range("P12").Value = 513

$ws.Range("D14").Value = 44453
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 12000
$ws.Range("P14").Value = 400
